$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.772.43"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "2.618.45"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.16"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.50"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.543"
$ws.Range("E8").Value = "  -2.45%  "
$ws.Range("D9").Value = "2.619.00"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("E10").Value = "  +6.44%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.52"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000187"
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("D16").Value = "3.088.05"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "67.648.27"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "2.617.41"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "371.89"
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.20"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.22"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("E22").Value = "  -13.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.79"
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -4.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.97"
$ws.Range("E25").Value = "  +10.02%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.87"
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000104"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "590.66"
$ws.Range("E30").Value = "  +1.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.74"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.16"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("E40").Value = "  +3.39%  "
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.69"
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.11"
$ws.Range("E44").Value = "  +4.66%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("D47").Value = "0.0₆0303"
$ws.Range("E47").Value = "  +5.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.35"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("E51").Value = "  -1.69%  "
